# Apply quarterly financial data updates to worksheet 'ZNH'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 10020600
$ws.Range("E8").Value = 9990700
$ws.Range("F8").Value = 8977000
$ws.Range("G8").Value = 9032500
$ws.Range("H8").Value = 8031800
$ws.Range("I8").Value = 8651300
$ws.Range("J8").Value = 7919000

$ws.Range("D9").Value = 7824000
$ws.Range("E9").Value = 7653400
$ws.Range("F9").Value = 6859700
$ws.Range("G9").Value = 6737100
$ws.Range("H9").Value = 5580100
$ws.Range("I9").Value = 6133200
$ws.Range("J9").Value = 5550700

$ws.Range("D10").Value = 2196600
$ws.Range("E10").Value = 2337300
$ws.Range("F10").Value = 2117400
$ws.Range("G10").Value = 2295500
$ws.Range("H10").Value = 2451700
$ws.Range("I10").Value = 2518100
$ws.Range("J10").Value = 2368300

$ws.Range("E14").Value = 31900
$ws.Range("I14").Value = 13400

$ws.Range("D15").Value = 1033200
$ws.Range("E15").Value = 1011400
$ws.Range("F15").Value = 965900
$ws.Range("G15").Value = 947000
$ws.Range("H15").Value = 945200
$ws.Range("I15").Value = 909200
$ws.Range("J15").Value = 870700

$ws.Range("D17").Value = 9327300
$ws.Range("E17").Value = 9295700
$ws.Range("F17").Value = 8297000
$ws.Range("G17").Value = 8193900
$ws.Range("H17").Value = 6985400
$ws.Range("I17").Value = 7545200
$ws.Range("J17").Value = 7030800

$ws.Range("D18").Value = 693400
$ws.Range("E18").Value = 695000
$ws.Range("F18").Value = 680000
$ws.Range("G18").Value = 838700
$ws.Range("H18").Value = 1046400
$ws.Range("I18").Value = 1106100
$ws.Range("J18").Value = 888200

$ws.Range("D20").Value = -13400
$ws.Range("E20").Value = 212100
$ws.Range("F20").Value = 134900
$ws.Range("G20").Value = -227800
$ws.Range("H20").Value = -161500
$ws.Range("I20").Value = -826100
$ws.Range("J20").Value = 64900

$ws.Range("D21").Value = 1713200
$ws.Range("E21").Value = 952600
$ws.Range("F21").Value = 1780800
$ws.Range("G21").Value = 612600
$ws.Range("H21").Value = 1830200
$ws.Range("I21").Value = 318500
$ws.Range("J21").Value = 1823800

$ws.Range("D22").Value = 228400
$ws.Range("E22").Value = 209400
$ws.Range("F22").Value = 195600
$ws.Range("G22").Value = 184300
$ws.Range("H22").Value = 174500
$ws.Range("I22").Value = 159400
$ws.Range("J22").Value = 165800

$ws.Range("D23").Value = 451600
$ws.Range("E23").Value = 697700
$ws.Range("F23").Value = 619300
$ws.Range("G23").Value = 426500
$ws.Range("H23").Value = 710400
$ws.Range("I23").Value = 120700
$ws.Range("J23").Value = 787300

$ws.Range("D24").Value = 101200
$ws.Range("E24").Value = 150500
$ws.Range("F24").Value = 142800
$ws.Range("G24").Value = 99000
$ws.Range("H24").Value = 162700
$ws.Range("I24").Value = 13400
$ws.Range("J24").Value = 179600

$ws.Range("D26").Value = 350400
$ws.Range("E26").Value = 547200
$ws.Range("F26").Value = 476500
$ws.Range("G26").Value = 327500
$ws.Range("H26").Value = 547800
$ws.Range("I26").Value = 107300
$ws.Range("J26").Value = 607700

$ws.Range("D27").Value = 310900
$ws.Range("E27").Value = 473300
$ws.Range("F27").Value = 411400
$ws.Range("G27").Value = 285800
$ws.Range("H27").Value = 462700
$ws.Range("I27").Value = 38100
$ws.Range("J27").Value = 516300

$ws.Range("D32").Value = 13400
$ws.Range("E32").Value = -212100
$ws.Range("F32").Value = -134900
$ws.Range("G32").Value = 227800
$ws.Range("H32").Value = 161500
$ws.Range("I32").Value = 826100
$ws.Range("J32").Value = -64900

$ws.Range("D33").Value = 310900
$ws.Range("E33").Value = 473300
$ws.Range("F33").Value = 411400
$ws.Range("G33").Value = 285800
$ws.Range("H33").Value = 462700
$ws.Range("I33").Value = 38100
$ws.Range("J33").Value = 516300

$ws.Range("D35").Value = 310900
$ws.Range("E35").Value = 473300
$ws.Range("F35").Value = 411400
$ws.Range("G35").Value = 285800
$ws.Range("H35").Value = 462700
$ws.Range("I35").Value = 38100
$ws.Range("J35").Value = 516300

$ws.Range("D41").Value = 805400
$ws.Range("E41").Value = 1013000
$ws.Range("F41").Value = 633400
$ws.Range("G41").Value = 616200
$ws.Range("H41").Value = 912900
$ws.Range("I41").Value = 676700
$ws.Range("J41").Value = 1425500

$ws.Range("D43").Value = 1558900
$ws.Range("E43").Value = 1184800
$ws.Range("F43").Value = 1117100
$ws.Range("G43").Value = 960800
$ws.Range("H43").Value = 1040200
$ws.Range("I43").Value = 984400
$ws.Range("J43").Value = 1276600

$ws.Range("D44").Value = 279000
$ws.Range("E44").Value = 240700
$ws.Range("F44").Value = 260500
$ws.Range("G44").Value = 235700
$ws.Range("H44").Value = 263900
$ws.Range("I44").Value = 238300
$ws.Range("J44").Value = 256500

$ws.Range("D45").Value = 238200
$ws.Range("E45").Value = 215600
$ws.Range("F45").Value = 243100
$ws.Range("G45").Value = 230000
$ws.Range("H45").Value = 228100
$ws.Range("I45").Value = 195000
$ws.Range("J45").Value = 145100

$ws.Range("D46").Value = 2881500
$ws.Range("E46").Value = 2654200
$ws.Range("F46").Value = 2254100
$ws.Range("G46").Value = 2042700
$ws.Range("H46").Value = 2445100
$ws.Range("I46").Value = 2094500
$ws.Range("J46").Value = 3103700

$ws.Range("D47").Value = 891400
$ws.Range("E47").Value = 803300
$ws.Range("F47").Value = 849900
$ws.Range("G47").Value = 807200
$ws.Range("H47").Value = 699800
$ws.Range("I47").Value = 689800
$ws.Range("J47").Value = 686700

$ws.Range("D48").Value = 29341800
$ws.Range("E48").Value = 28073100
$ws.Range("F48").Value = 26376800
$ws.Range("G48").Value = 26069100
$ws.Range("H48").Value = 24389800
$ws.Range("I48").Value = 24087400
$ws.Range("J48").Value = 23502400

$ws.Range("D49").Value = 35200
$ws.Range("E49").Value = 35200
$ws.Range("F49").Value = 27000
$ws.Range("G49").Value = 27000

$ws.Range("D52").Value = 904600
$ws.Range("E52").Value = 894200
$ws.Range("F52").Value = 767900
$ws.Range("G52").Value = 801600
$ws.Range("H52").Value = 806600
$ws.Range("I52").Value = 730900
$ws.Range("J52").Value = 634300

$ws.Range("D54").Value = 34054500
$ws.Range("E54").Value = 32459900
$ws.Range("F54").Value = 30275600
$ws.Range("G54").Value = 29747600
$ws.Range("H54").Value = 28341300
$ws.Range("I54").Value = 27602600
$ws.Range("J54").Value = 27927100

$ws.Range("D57").Value = 300100
$ws.Range("E57").Value = 315400
$ws.Range("F57").Value = 290400
$ws.Range("G57").Value = 282400
$ws.Range("H57").Value = 263100
$ws.Range("I57").Value = 371000
$ws.Range("J57").Value = 310300

$ws.Range("D58").Value = 7243300
$ws.Range("E58").Value = 5329300
$ws.Range("F58").Value = 4447800
$ws.Range("G58").Value = 5259800
$ws.Range("H58").Value = 5630500
$ws.Range("I58").Value = 5404800
$ws.Range("J58").Value = 4074700

$ws.Range("D59").Value = 4632900
$ws.Range("E59").Value = 4681300
$ws.Range("F59").Value = 4490900
$ws.Range("G59").Value = 4539600
$ws.Range("H59").Value = 4128800
$ws.Range("I59").Value = 3950200
$ws.Range("J59").Value = 3701800

$ws.Range("D60").Value = 12176300
$ws.Range("E60").Value = 10325900
$ws.Range("F60").Value = 9229200
$ws.Range("G60").Value = 10081800
$ws.Range("H60").Value = 10022400
$ws.Range("I60").Value = 9726000
$ws.Range("J60").Value = 8086900

$ws.Range("D61").Value = 11228000
$ws.Range("E61").Value = 11917600
$ws.Range("F61").Value = 11625100
$ws.Range("G61").Value = 10727800
$ws.Range("H61").Value = 9690900
$ws.Range("I61").Value = 9690000
$ws.Range("J61").Value = 11690700

$ws.Range("D62").Value = 1026800
$ws.Range("E62").Value = 934400
$ws.Range("F62").Value = 879500
$ws.Range("G62").Value = 779000
$ws.Range("H62").Value = 832300
$ws.Range("I62").Value = 821900
$ws.Range("J62").Value = 829000

$ws.Range("D66").Value = 26365600
$ws.Range("E66").Value = 25048900
$ws.Range("F66").Value = 23544500
$ws.Range("G66").Value = 23298300
$ws.Range("H66").Value = 22215000
$ws.Range("I66").Value = 21808000
$ws.Range("J66").Value = 22162500

$ws.Range("D72").Value = 3865000
$ws.Range("E72").Value = 3619900
$ws.Range("F72").Value = 3146400
$ws.Range("G72").Value = 2864000
$ws.Range("H72").Value = 2578500
$ws.Range("I72").Value = 2232200
$ws.Range("J72").Value = 2201700

$ws.Range("D76").Value = 7688800
$ws.Range("E76").Value = 7411000
$ws.Range("F76").Value = 6731100
$ws.Range("G76").Value = 6449300
$ws.Range("H76").Value = 6126200
$ws.Range("I76").Value = 5794700
$ws.Range("J76").Value = 5764500

$ws.Range("D81").Value = 310900
$ws.Range("E81").Value = 473300
$ws.Range("F81").Value = 411400
$ws.Range("G81").Value = 285800
$ws.Range("H81").Value = 462700
$ws.Range("I81").Value = 38100
$ws.Range("J81").Value = 516300

$ws.Range("D89").Value = 1092700
$ws.Range("E89").Value = 1575100
$ws.Range("F89").Value = 1056500
$ws.Range("G89").Value = 2322800
$ws.Range("H89").Value = 1204100
$ws.Range("I89").Value = 2055300
$ws.Range("J89").Value = 1467000

$ws.Range("D91").Value = -1403400
$ws.Range("E91").Value = -1248000
$ws.Range("F91").Value = -806900
$ws.Range("G91").Value = -2298700
$ws.Range("H91").Value = -516200
$ws.Range("I91").Value = -1157600
$ws.Range("J91").Value = -644000

$ws.Range("D94").Value = -1003400
$ws.Range("E94").Value = -787800
$ws.Range("F94").Value = -434500
$ws.Range("G94").Value = -2154500
$ws.Range("H94").Value = -183000
$ws.Range("I94").Value = -454700
$ws.Range("J94").Value = -573900

$ws.Range("E96").Value = -145700
$ws.Range("G96").Value = -49700
$ws.Range("H96").Value = -66800

$ws.Range("D100").Value = -297400
$ws.Range("E100").Value = -404700
$ws.Range("F100").Value = -603900
$ws.Range("G100").Value = -467600
$ws.Range("H100").Value = -787800
$ws.Range("I100").Value = -2357000
$ws.Range("J100").Value = -1753200

$ws.Range("E101").Value = -3000
$ws.Range("G101").Value = 2700
$ws.Range("I101").Value = 7700
$ws.Range("J101").Value = -2100

$ws.Range("D102").Value = -207600
$ws.Range("E102").Value = 379600
$ws.Range("F102").Value = 17200
$ws.Range("G102").Value = -296700
$ws.Range("H102").Value = 236100
$ws.Range("I102").Value = -748700
$ws.Range("J102").Value = -862100
